$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 166670300
$ws.Range("I69").Value = 250002940
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 750008820
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -750007946
$ws.Range("N69").Value = -16748
$ws.Range("H72").Value = 166670300
$ws.Range("I72").Value = 250002940
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 2250026460
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -2250022092
$ws.Range("N72").Value = -53736
$ws.Range("H87").Value = 81687.5
$ws.Range("J87").Value = 81687.5
$ws.Range("L87").Value = 81687.5
$ws.Range("N87").Value = -84183.5
$ws.Range("H90").Value = 81687.5
$ws.Range("J90").Value = 81687.5
$ws.Range("L90").Value = 245062.5
$ws.Range("N90").Value = -257542.5
$ws.Range("H132").Value = 176264.88
$ws.Range("I132").Value = 284607.75
$ws.Range("J132").Value = 18675.227
$ws.Range("K132").Value = 853823.25
$ws.Range("L132").Value = 56025.681
$ws.Range("M132").Value = -851293.25
$ws.Range("N132").Value = -61085.681
$ws.Range("H135").Value = 4003.6897
$ws.Range("J135").Value = 5174.421
$ws.Range("L135").Value = 46569.789
$ws.Range("N135").Value = -51639.789
$ws.Range("H138").Value = 3807.42
$ws.Range("I138").Value = 1997.4073
$ws.Range("J138").Value = 4476.877
$ws.Range("K138").Value = 5992.2219
$ws.Range("L138").Value = 13430.631
$ws.Range("M138").Value = -852.2219000000005
$ws.Range("N138").Value = -23710.631
$ws.Range("H141").Value = 5485.9443
$ws.Range("I141").Value = 5596.647
$ws.Range("K141").Value = 16789.941
$ws.Range("M141").Value = -11609.941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17201.037
$ws.Range("I32").Value = 17001.229
$ws.Range("J32").Value = 18571.143
$ws.Range("K32").Value = 17001.229
$ws.Range("L32").Value = 18571.143
$ws.Range("M32").Value = -16714.229
$ws.Range("N32").Value = -19145.143
$ws.Range("H35").Value = 3730
$ws.Range("I35").Value = 3730
$ws.Range("K35").Value = 3730
$ws.Range("M35").Value = -3324
$ws.Range("H74").Value = 10001558
$ws.Range("I74").Value = 14706988
$ws.Range("K74").Value = 14706988
$ws.Range("M74").Value = -14706114
$ws.Range("H77").Value = 10001558
$ws.Range("I77").Value = 14706988
$ws.Range("K77").Value = 73534940
$ws.Range("M77").Value = -73530572
$ws.Range("H80").Value = 77247
$ws.Range("J80").Value = 77247
$ws.Range("L80").Value = 77247
$ws.Range("N80").Value = -79243
$ws.Range("H83").Value = 77247
$ws.Range("J83").Value = 77247
$ws.Range("L83").Value = 231741
$ws.Range("N83").Value = -241725
$ws.Range("H132").Value = 19170.564
$ws.Range("J132").Value = 9472.532999999999
$ws.Range("L132").Value = 28417.599
$ws.Range("N132").Value = -33477.599

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1851.5555
$ws.Range("I134").Value = 1850.5
$ws.Range("K134").Value = 5551.5
$ws.Range("M134").Value = -3016.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 568.1579
$ws.Range("I22").Value = 226.5
$ws.Range("J22").Value = 1153.8572
$ws.Range("K22").Value = 226.5
$ws.Range("L22").Value = 1153.8572
$ws.Range("M22").Value = 123.5
$ws.Range("N22").Value = -1853.8572
$ws.Range("H31").Value = 14088382
$ws.Range("I31").Value = 31251926
$ws.Range("J31").Value = 5474.231
$ws.Range("K31").Value = 31251926
$ws.Range("L31").Value = 5474.231
$ws.Range("M31").Value = -31251631
$ws.Range("N31").Value = -6064.231
$ws.Range("H34").Value = 14088382
$ws.Range("I34").Value = 31251926
$ws.Range("J34").Value = 5474.231
$ws.Range("K34").Value = 31251926
$ws.Range("L34").Value = 5474.231
$ws.Range("M34").Value = -31251724
$ws.Range("N34").Value = -5878.231
$ws.Range("H74").Value = 63250.668
$ws.Range("J74").Value = 64100.8
$ws.Range("L74").Value = 64100.8
$ws.Range("N74").Value = -65848.8
$ws.Range("H77").Value = 63250.668
$ws.Range("J77").Value = 64100.8
$ws.Range("L77").Value = 192302.4
$ws.Range("N77").Value = -201038.4
$ws.Range("H105").Value = 1749764.8
$ws.Range("I105").Value = 2273961
$ws.Range("K105").Value = 2273961
$ws.Range("M105").Value = -2272214
$ws.Range("H132").Value = 21510372
$ws.Range("I132").Value = 24693172
$ws.Range("K132").Value = 74079516
$ws.Range("M132").Value = -74076986
$ws.Range("H133").Value = 76609.91
$ws.Range("J133").Value = 75847
$ws.Range("L133").Value = 75847
$ws.Range("N133").Value = -80907

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 344.42856
$ws.Range("I2").Value = 349.66666
$ws.Range("K2").Value = 2097.99996
$ws.Range("M2").Value = -1984.99996
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H129").Value = 1094.25
$ws.Range("J129").Value = 2299
$ws.Range("L129").Value = 6897
$ws.Range("N129").Value = -16897

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I122").Value = 459491.72
$ws.Range("K122").Value = 1378475.16
$ws.Range("M122").Value = -1376025.16
$ws.Range("H132").Value = 5520.174
$ws.Range("I132").Value = 3902.6428
$ws.Range("K132").Value = 11707.9284
$ws.Range("M132").Value = -9177.928400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1423954.8
$ws.Range("I68").Value = 2842551.5
$ws.Range("J68").Value = 5358.125
$ws.Range("K68").Value = 2842551.5
$ws.Range("L68").Value = 5358.125
$ws.Range("M68").Value = -2841802.5
$ws.Range("N68").Value = -6856.125
$ws.Range("H71").Value = 1423954.8
$ws.Range("I71").Value = 2842551.5
$ws.Range("J71").Value = 5358.125
$ws.Range("K71").Value = 14212757.5
$ws.Range("L71").Value = 26790.625
$ws.Range("M71").Value = -14209013.5
$ws.Range("N71").Value = -34278.625
$ws.Range("H80").Value = 99950
$ws.Range("J80").Value = 99950
$ws.Range("L80").Value = 99950
$ws.Range("N80").Value = -102196
$ws.Range("H83").Value = 99950
$ws.Range("J83").Value = 99950
$ws.Range("L83").Value = 299850
$ws.Range("N83").Value = -311082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 7000
$ws.Range("J8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("N8").Value = -7280
$ws.Range("H80").Value = 101633.336
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 101633.336
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 101633.336
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -103629.336
$ws.Range("H83").Value = 101633.336
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 101633.336
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 304900.008
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -314884.008
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H132").Value = 9617252
$ws.Range("I132").Value = 1053.8235
$ws.Range("K132").Value = 3161.4705
$ws.Range("M132").Value = -631.4704999999999
